# NIT-9001720691.xlsx - "Actualiza base de datos EC y agrega parte 1 de
# nuevos estado de cuenta" (Update EC database and add part 1 of new
# account statement).
#
# The old worksheet carried two workers' overdue-period tables:
#   - TANIA MILENA CHAMORRO VANEGAS (CC 52793304), 8 periods (rows 16-23)
#   - ADONYS JUNIOR CHAVEZ ROSSO (CC 1143350645), 3 periods (rows 24-26)
# The update drops TANIA's block entirely and keeps only ADONYS' 3 rows
# (now re-ordered oldest-period-first), along with updated summary
# figures at the top of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove TANIA's 8 data rows (old rows 16-23). This shifts the ADONYS
# rows (old 24-26) up to become the new rows 16-18, and shifts the
# signature-block rows (old 31-32) up to become the new rows 23-24 -
# exactly matching the target layout.
$ws.Rows("16:23").Delete()

# --- Updated summary figures (top of the statement) ---
$ws.Range("E11").Value = 120000   # VALOR MORA
$ws.Range("C13").Value = 1        # Cant. Trabajadores
$ws.Range("F13").Value = 3        # Cant. Periodos

# --- Re-order ADONYS' periods to ascending order (1909, 1910, 1911) ---
$ws.Range("E16").Value = "1909"
$ws.Range("E17").Value = "1910"
$ws.Range("E18").Value = "1911"

# --- Column D ("Nombre Trabajador") shrinks now that the longer TANIA
# name is gone and the best-fit width should reflect the ADONYS name. ---
$ws.Columns("D:D").AutoFit()
